# Update tab names in all BOMs, fix bi-color LED naming.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab from "kNOT" to "BOM"
$ws.Name = "BOM"

# The "R13, R14, R15" / "R10, R11, R12" / "D1, D2, D3" / "LD_A, LD_B, LD_C" /
# "J_IN_A, ..." / "J_B1, J_B3" reference-designator cells in column D were
# previously (incorrectly) styled with an explicit "plain font" style that is
# no longer needed -- reset them back to the workbook's Normal style.
$refDesCells = @("D3", "D4", "D10", "D15", "D16", "D20")
foreach ($addr in $refDesCells) {
    $ws.Range($addr).Style = "Normal"
}

# Move the active cell / selection to D34 (matches latest saved view state).
$ws.Range("D34").Select() | Out-Null
